# Update the "Listing" table to match the new schema:
#  - rename header "price" (F1) to "Rent"
#  - update the F-column values to the new rent figures
#  - move the active selection to F2 (matches the author's saved cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Rent"

$ws.Range("F2").Value = 2456
$ws.Range("F3").Value = 3789
$ws.Range("F4").Value = 2578

$ws.Range("F2").Select()
